$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text-typed cell value without Excel auto-converting
# date-like strings (e.g. "2025-08-29") into date serial numbers, and
# without leaving a stray NumberFormat style applied to the cell.
function Set-TextCell {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

function Set-NumberCell {
    param($cell, $value)
    $cell.Value = $value
}

# Row 3
Set-NumberCell $ws.Range("B3") 79244

# Row 4
Set-NumberCell $ws.Range("A4") 130134267
Set-NumberCell $ws.Range("B4") 99014
Set-TextCell $ws.Range("D4") "VU"
Set-NumberCell $ws.Range("E4") 220787
Set-TextCell $ws.Range("F4") "Knärot"
Set-TextCell $ws.Range("G4") "Goodyera repens"
Set-TextCell $ws.Range("H4") "(L.) R. Br."
Set-NumberCell $ws.Range("Q4") 750666
Set-NumberCell $ws.Range("R4") 7111136
Set-TextCell $ws.Range("Y4") "2025-08-29"
Set-TextCell $ws.Range("AA4") "2025-08-29"
Set-TextCell $ws.Range("AX4") "Lisa Sandberg"

# Row 5
Set-NumberCell $ws.Range("A5") 130134356
Set-NumberCell $ws.Range("B5") 93096
Set-TextCell $ws.Range("D5") "LC"
Set-NumberCell $ws.Range("E5") 4364
Set-TextCell $ws.Range("F5") "Dropptaggsvamp"
Set-TextCell $ws.Range("G5") "Hydnellum ferrugineum"
Set-TextCell $ws.Range("H5") "(Fr.:Fr.) P. Karst."
Set-NumberCell $ws.Range("Q5") 750732
Set-NumberCell $ws.Range("R5") 7111314
Set-TextCell $ws.Range("Y5") "2025-10-07"
Set-TextCell $ws.Range("AA5") "2025-10-07"
Set-TextCell $ws.Range("AX5") "Daniel Lussetti"

# Row 6
Set-NumberCell $ws.Range("B6") 79244

# Row 7
Set-NumberCell $ws.Range("B7") 93096

# Row 8
Set-NumberCell $ws.Range("B8") 79244

# Row 9
Set-NumberCell $ws.Range("B9") 79244

# Row 10
Set-NumberCell $ws.Range("A10") 130134260
Set-NumberCell $ws.Range("B10") 79001
Set-NumberCell $ws.Range("E10") 6446
Set-TextCell $ws.Range("F10") "Kolflarnlav"
Set-TextCell $ws.Range("G10") "Carbonicola anthracophila"
Set-TextCell $ws.Range("H10") "(Nyl.) Bendiksby & Timdal"
Set-NumberCell $ws.Range("Q10") 750719
Set-NumberCell $ws.Range("R10") 7111349

# Row 11
Set-NumberCell $ws.Range("A11") 130134263
Set-NumberCell $ws.Range("B11") 79244
Set-NumberCell $ws.Range("E11") 6425
Set-TextCell $ws.Range("F11") "Garnlav"
Set-TextCell $ws.Range("G11") "Alectoria sarmentosa"
Set-TextCell $ws.Range("H11") "(Ach.) Ach."
Set-NumberCell $ws.Range("Q11") 750665
Set-NumberCell $ws.Range("R11") 7111235

# Row 12
Set-NumberCell $ws.Range("B12") 93108

# Row 13
Set-NumberCell $ws.Range("A13") 130134271
Set-NumberCell $ws.Range("B13") 79244
Set-TextCell $ws.Range("D13") "NT"
Set-NumberCell $ws.Range("E13") 6425
Set-TextCell $ws.Range("F13") "Garnlav"
Set-TextCell $ws.Range("G13") "Alectoria sarmentosa"
Set-TextCell $ws.Range("H13") "(Ach.) Ach."
Set-NumberCell $ws.Range("Q13") 750620
Set-NumberCell $ws.Range("R13") 7111033

# Row 14
Set-NumberCell $ws.Range("A14") 130134238
Set-NumberCell $ws.Range("B14") 99014
Set-TextCell $ws.Range("D14") "VU"
Set-NumberCell $ws.Range("E14") 220787
Set-TextCell $ws.Range("F14") "Knärot"
Set-TextCell $ws.Range("G14") "Goodyera repens"
Set-TextCell $ws.Range("H14") "(L.) R. Br."
Set-NumberCell $ws.Range("Q14") 750467
Set-NumberCell $ws.Range("R14") 7110939

# Row 15
Set-NumberCell $ws.Range("B15") 79244

# Row 16
Set-NumberCell $ws.Range("B16") 78647

# Row 17
Set-NumberCell $ws.Range("B17") 79244

# Row 19
Set-NumberCell $ws.Range("B19") 93134

# Row 20
Set-NumberCell $ws.Range("B20") 93134
